# Muharrem - TestCase added.
#
# Adds a new Cucumber-style test case ("SHADE KADIN URUNLERI" / Blouses
# discount flow) to the bottom of the LCW test-case tracker sheet, reusing
# the same row-group layout already used for the other test cases above
# (header cell block, numbered "Steps" list, "Expected result" column,
# "Test Result" = successful/PASSED, then a Cucumber/Gherkin restatement
# of the same steps as Given/When/And/Then lines).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 155: replace the placeholder test-case summary / first step ---
$ws.Range("C155").Value = "The user should be able to select any of the Blouses from the SHADE WOMEN'S PRODUCTS category and go to the payment page. They should be able to see if there is a discount."
$ws.Range("E155").Value = "1. Open website / go to home page "
$ws.Range("G155").Value = "LCW home page should be open /display"
$ws.Range("H155").Value = "successful"
$ws.Range("I155").Value = "PASSED"

# Row 155 now carries a long wrapped description, so it needs a much
# taller row (matches the rest of the "Steps"/"Expected result" rows that
# already use wrapped text on this sheet).
$ws.Range("A155").EntireRow.RowHeight = 90.6

# --- Rows 156-165: remaining numbered steps + expected results ---
$ws.Range("E156").Value = "2. Verify that you are on the home page"
$ws.Range("G156").Value = "The user must be verified to be on the homepage"
$ws.Range("H156").Value = "successful"

$ws.Range("E157").Value = "3. Click on `"SHADE KADIN ÜRÜNLERİ`" title opportunity"
$ws.Range("G157").Value = "User should be able to click on `"SHADE KADIN ÜRÜNLERİ`" title"
$ws.Range("H157").Value = "successful"

$ws.Range("E158").Value = "4. Filter products by `"Bluz`" and choose a random product"
$ws.Range("G158").Value = "User should be able to choose `"Bluz`" filter "
$ws.Range("H158").Value = "successful"

$ws.Range("E159").Value = "5. If available, choose random size and color and add to cart"
$ws.Range("G159").Value = "User should be able if available choose random size and color and add to cart"
$ws.Range("H159").Value = "successful"

$ws.Range("E160").Value = "6. Click on the addToCart"
$ws.Range("G160").Value = "User should be able to click the any bot "
$ws.Range("H160").Value = "successful"

$ws.Range("E161").Value = "7. Click on the Cart"
$ws.Range("G161").Value = "User should be able to click the any size"
$ws.Range("H161").Value = "successful"

$ws.Range("E162").Value = "8. Product price and discount rate"
$ws.Range("G162").Value = "User should be able to click the add to cart button"
$ws.Range("H162").Value = "successful"

$ws.Range("E163").Value = "9. Click on the Payment"
$ws.Range("G163").Value = "User should be able to click the Cart "
$ws.Range("H163").Value = "successful"

$ws.Range("E164").Value = "10. Fill out the form on the payment page"
$ws.Range("G164").Value = "User should be able to save the product price"
$ws.Range("H164").Value = "successful"

$ws.Range("E165").Value = "11. Click on the save button"
$ws.Range("G165").Value = "User should be able to click the Save Button "
$ws.Range("H165").Value = "successful"

# --- Rows 167-168: "Cucumber type test case" / "Gherkin language" labels ---
# These reuse the bold label style already used elsewhere on the sheet for
# this sub-heading pair, which this workbook doesn't have yet -> Excel
# creates the new bold font + cell style automatically the first time it's
# applied here.
$ws.Range("E167").Value = "Cucumber type test case"
$ws.Range("E167").Font.Bold = $true

$ws.Range("E168").Value = "Gherkin language"
$ws.Range("E168").Font.Bold = $true

# --- Rows 169-180: Gherkin restatement of the same steps ---
$ws.Range("E169").Value = "GIVEN User open app / go to home page "
$ws.Range("E170").Value = "WHEN Verify that you are on the home page"
$ws.Range("E171").Value = "And Click on MARKALARA ÖZEL title from the main page"
$ws.Range("E172").Value = "And Click on SHADE KADIN ÜRÜNLERİ title opportunity"
$ws.Range("E173").Value = "And Filter products by Bluz and choose a random product"
$ws.Range("E174").Value = "And If available, choose random size and color and add to cart"
$ws.Range("E175").Value = "And Click on the addToCart"
$ws.Range("E176").Value = "And Click on the Cart"
$ws.Range("E177").Value = "And Product price and discount rate"
$ws.Range("E178").Value = "And Click on the Payment"
$ws.Range("E179").Value = "And Fill out the form on the payment page"
$ws.Range("E180").Value = "Then Click on the save button"

# --- View state: leave the selection on the last filled cell, scrolled
# down to the new test case, matching where the author was working ---
$ws.Range("J155").Select()
$excel.ActiveWindow.Zoom = 50
